$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 - fill in actual clock times
$ws.Range("C53").Value = 0.26250000000000001
$ws.Range("D53").Value = 0.50555555555555554
$ws.Range("E53").Value = 0.52986111111111112
$ws.Range("F53").Value = 0.67013888888888884

# Row 54 - fill in actual clock times
$ws.Range("C54").Value = 0.28333333333333333
$ws.Range("D54").Value = 0.51527777777777783
$ws.Range("E54").Value = 0.55347222222222225
$ws.Range("F54").Value = 0.64583333333333337

# Row 63 - fill in actual clock times (only clock-in / clock-out, no lunch break)
$ws.Range("C63").Value = 0.28541666666666665
$ws.Range("F63").Value = 0.61458333333333337

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 50
$ws.Range("G61").Select()
